$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.052343858745701
$ws.Cells.Item(2, 3).Value = 0.1875371153331002
$ws.Cells.Item(2, 4).Value = 0.07843894618501679
$ws.Cells.Item(2, 5).Value = 0.0894042697841797
$ws.Cells.Item(2, 7).Value = 0.002449270318185254
$ws.Cells.Item(2, 9).Value = 0.8132819134243903
$ws.Cells.Item(2, 12).Value = 0.2075585828297832
$ws.Cells.Item(2, 13).Value = 0.238389857254262
$ws.Cells.Item(2, 14).Value = 1.409308619173
$ws.Cells.Item(2, 15).Value = 3.45736719679735
$ws.Cells.Item(3, 2).Value = 0.9657966250594541
$ws.Cells.Item(3, 3).Value = 0.1701085491547474
$ws.Cells.Item(3, 4).Value = 0.0711607546129045
$ws.Cells.Item(3, 5).Value = 0.08997966645856592
$ws.Cells.Item(3, 7).Value = 0.002452490693154787
$ws.Cells.Item(3, 9).Value = 0.818311894595162
$ws.Cells.Item(3, 12).Value = 0.2048095850056555
$ws.Cells.Item(3, 13).Value = 0.2246725768106259
$ws.Cells.Item(3, 14).Value = 1.424864717837735
$ws.Cells.Item(3, 15).Value = 3.447839505702461
$ws.Cells.Item(4, 2).Value = 0.9129031115281236
$ws.Cells.Item(4, 3).Value = 0.1593274632390091
$ws.Cells.Item(4, 4).Value = 0.06672771432938873
$ws.Cells.Item(4, 5).Value = 0.09035340718596996
$ws.Cells.Item(4, 7).Value = 0.002454574854748242
$ws.Cells.Item(4, 9).Value = 0.8218745917599932
$ws.Cells.Item(4, 12).Value = 0.2032254755952252
$ws.Cells.Item(4, 13).Value = 0.216337941850739
$ws.Cells.Item(4, 14).Value = 1.434954338134013
$ws.Cells.Item(4, 15).Value = 3.444215274260131
$ws.Cells.Item(5, 2).Value = 0.8914117087055331
$ws.Cells.Item(5, 3).Value = 0.1549140874368788
$ws.Cells.Item(5, 4).Value = 0.06493018536927764
$ws.Cells.Item(5, 5).Value = 0.09051085975938533
$ws.Cells.Item(5, 7).Value = 0.002455451114964365
$ws.Cells.Item(5, 9).Value = 0.8234455709790502
$ws.Cells.Item(5, 12).Value = 0.202606080289975
$ws.Cells.Item(5, 13).Value = 0.2129637644178075
$ws.Cells.Item(5, 14).Value = 1.439201258419725
$ws.Cells.Item(5, 15).Value = 3.443297439859009
$ws.Cells.Item(6, 2).Value = 0.8878469221724288
$ws.Cells.Item(6, 3).Value = 0.1541800447167816
$ws.Cells.Item(6, 4).Value = 0.06463224797890632
$ws.Cells.Item(6, 5).Value = 0.0905373159680094
$ws.Cells.Item(6, 7).Value = 0.002455598247282425
$ws.Cells.Item(6, 9).Value = 0.8237136241907059
$ws.Cells.Item(6, 12).Value = 0.2025048104014431
$ws.Cells.Item(6, 13).Value = 0.2124048338831344
$ws.Cells.Item(6, 14).Value = 1.439914629408907
$ws.Cells.Item(6, 15).Value = 3.443178782279858
$ws.Cells.Item(7, 2).Value = 0.9126130139177633
$ws.Cells.Item(7, 3).Value = 0.1592680237508262
$ws.Cells.Item(7, 4).Value = 0.06670343594225869
$ws.Cells.Item(7, 5).Value = 0.09035550978111961
$ws.Cells.Item(7, 7).Value = 0.002454586563097148
$ws.Cells.Item(7, 9).Value = 0.8218952962409247
$ws.Cells.Item(7, 12).Value = 0.2032170163051532
$ws.Cells.Item(7, 13).Value = 0.2162923461849644
$ws.Cells.Item(7, 14).Value = 1.435011065682055
$ws.Cells.Item(7, 15).Value = 3.444200633215303
$ws.Cells.Item(8, 2).Value = 1.022451750871767
$ws.Cells.Item(8, 3).Value = 0.18154437277542
$ws.Cells.Item(8, 4).Value = 0.07592197484440533
$ws.Cells.Item(8, 5).Value = 0.08959842878285162
$ws.Cells.Item(8, 7).Value = 0.002450358578654055
$ws.Cells.Item(8, 9).Value = 0.8149177509350807
$ws.Cells.Item(8, 12).Value = 0.2065892126129256
$ws.Cells.Item(8, 13).Value = 0.2336420058127686
$ws.Cells.Item(8, 14).Value = 1.414560680339413
$ws.Cells.Item(8, 15).Value = 3.453619689958117
$ws.Cells.Item(9, 2).Value = 1.23976827902294
$ws.Cells.Item(9, 3).Value = 0.2245924775766071
$ws.Cells.Item(9, 4).Value = 0.09428590218382737
$ws.Cells.Item(9, 5).Value = 0.08827559113032768
$ws.Cells.Item(9, 7).Value = 0.002442911483852271
$ws.Cells.Item(9, 9).Value = 0.8050032482231515
$ws.Cells.Item(9, 12).Value = 0.2140242305482971
$ws.Cells.Item(9, 13).Value = 0.2683560515254442
$ws.Cells.Item(9, 14).Value = 1.378726612629297
$ws.Cells.Item(9, 15).Value = 3.48978588925695
$ws.Cells.Item(10, 2).Value = 1.400572137335416
$ws.Cells.Item(10, 3).Value = 0.2558326173441969
$ws.Cells.Item(10, 4).Value = 0.1079571141624029
$ws.Cells.Item(10, 5).Value = 0.0874017339444213
$ws.Cells.Item(10, 7).Value = 0.002437949367938359
$ws.Cells.Item(10, 9).Value = 0.8000245000666339
$ws.Cells.Item(10, 12).Value = 0.2199870141572831
$ws.Cells.Item(10, 13).Value = 0.2942774209108592
$ws.Cells.Item(10, 14).Value = 1.355000492088845
$ws.Cells.Item(10, 15).Value = 3.527201898375466
$ws.Cells.Item(11, 2).Value = 1.473968320568133
$ws.Cells.Item(11, 3).Value = 0.2699607024801196
$ws.Cells.Item(11, 4).Value = 0.1142164934257295
$ws.Cells.Item(11, 5).Value = 0.08702535159881075
$ws.Cells.Item(11, 7).Value = 0.002435801439776066
$ws.Cells.Item(11, 9).Value = 0.7982619491077116
$ws.Cells.Item(11, 12).Value = 0.2228081314950003
$ws.Cells.Item(11, 13).Value = 0.306159479490546
$ws.Cells.Item(11, 14).Value = 1.344771455226653
$ws.Cells.Item(11, 15).Value = 3.54659113966764
$ws.Cells.Item(12, 2).Value = 1.501796040894703
$ws.Cells.Item(12, 3).Value = 0.2752986254452594
$ws.Cells.Item(12, 4).Value = 0.1165926033560964
$ws.Cells.Item(12, 5).Value = 0.08688585582093666
$ws.Cells.Item(12, 7).Value = 0.002435003716960297
$ws.Cells.Item(12, 9).Value = 0.7976668874563018
$ws.Cells.Item(12, 12).Value = 0.2238920043638899
$ws.Cells.Item(12, 13).Value = 0.3106717621412329
$ws.Cells.Item(12, 14).Value = 1.340979119184045
$ws.Cells.Item(12, 15).Value = 3.554274805952417
$ws.Cells.Item(13, 2).Value = 1.495801341636366
$ws.Cells.Item(13, 3).Value = 0.2741495470803272
$ws.Cells.Item(13, 4).Value = 0.1160806063456903
$ws.Cells.Item(13, 5).Value = 0.08691576400526291
$ws.Cells.Item(13, 7).Value = 0.002435174825994691
$ws.Cells.Item(13, 9).Value = 0.7977918232447365
$ws.Cells.Item(13, 12).Value = 0.2236578808274032
$ws.Cells.Item(13, 13).Value = 0.3096993946906821
$ws.Cells.Item(13, 14).Value = 1.341792255143837
$ws.Cells.Item(13, 15).Value = 3.552604795274505
$ws.Cells.Item(14, 2).Value = 1.47625704546823
$ws.Cells.Item(14, 3).Value = 0.2704000992728197
$ws.Cells.Item(14, 4).Value = 0.1144118606468396
$ws.Cells.Item(14, 5).Value = 0.08701381447492695
$ws.Cells.Item(14, 7).Value = 0.00243573549738374
$ws.Cells.Item(14, 9).Value = 0.7982115416175688
$ws.Cells.Item(14, 12).Value = 0.2228969904229388
$ws.Cells.Item(14, 13).Value = 0.3065304520976326
$ws.Cells.Item(14, 14).Value = 1.344457830154152
$ws.Cells.Item(14, 15).Value = 3.547216433182143
$ws.Cells.Item(15, 2).Value = 1.46429001664427
$ws.Cells.Item(15, 3).Value = 0.2681018782662647
$ws.Cells.Item(15, 4).Value = 0.1133904643070593
$ws.Cells.Item(15, 5).Value = 0.08707426785007299
$ws.Cells.Item(15, 7).Value = 0.002436080960799961
$ws.Cells.Item(15, 9).Value = 0.7984780613842233
$ws.Cells.Item(15, 12).Value = 0.222432950414003
$ws.Cells.Item(15, 13).Value = 0.3045910455133054
$ws.Cells.Item(15, 14).Value = 1.346101145535155
$ws.Cells.Item(15, 15).Value = 3.543960387291406
$ws.Cells.Item(16, 2).Value = 1.395780369656336
$ws.Cells.Item(16, 3).Value = 0.2549076298735145
$ws.Cells.Item(16, 4).Value = 0.1075488618027123
$ws.Cells.Item(16, 5).Value = 0.08742675613101725
$ws.Cells.Item(16, 7).Value = 0.002438091934391408
$ws.Cells.Item(16, 9).Value = 0.800149807169646
$ws.Cells.Item(16, 12).Value = 0.2198048295611414
$ws.Cells.Item(16, 13).Value = 0.2935027020806871
$ws.Cells.Item(16, 14).Value = 1.355680339933517
$ws.Cells.Item(16, 15).Value = 3.525982494060401
$ws.Cells.Item(17, 2).Value = 1.353814028878162
$ws.Cells.Item(17, 3).Value = 0.2467919947150108
$ws.Cells.Item(17, 4).Value = 0.1039755630791888
$ws.Cells.Item(17, 5).Value = 0.08764840507490179
$ws.Cells.Item(17, 7).Value = 0.002439353558315445
$ws.Cells.Item(17, 9).Value = 0.8013041263103702
$ws.Cells.Item(17, 12).Value = 0.218220352038756
$ws.Cells.Item(17, 13).Value = 0.2867233618690292
$ws.Cells.Item(17, 14).Value = 1.36170138631368
$ws.Cells.Item(17, 15).Value = 3.515560834765665
$ws.Cells.Item(18, 2).Value = 1.329699294918953
$ws.Cells.Item(18, 3).Value = 0.2421162703988955
$ws.Cells.Item(18, 4).Value = 0.1019240855731596
$ws.Cells.Item(18, 5).Value = 0.08777788191476388
$ws.Cells.Item(18, 7).Value = 0.002440089509364969
$ws.Cells.Item(18, 9).Value = 0.8020153307871709
$ws.Cells.Item(18, 12).Value = 0.2173192285446532
$ws.Cells.Item(18, 13).Value = 0.2828325746390448
$ws.Cells.Item(18, 14).Value = 1.365217624159161
$ws.Cells.Item(18, 15).Value = 3.509789459964395
$ws.Cells.Item(19, 2).Value = 1.321538479350636
$ws.Cells.Item(19, 3).Value = 0.240531809799279
$ws.Cells.Item(19, 4).Value = 0.1012301405387603
$ws.Cells.Item(19, 5).Value = 0.08782206261073133
$ws.Cells.Item(19, 7).Value = 0.00244034046060899
$ws.Cells.Item(19, 9).Value = 0.8022642467107417
$ws.Cells.Item(19, 12).Value = 0.217015881233749
$ws.Cells.Item(19, 13).Value = 0.2815166880141504
$ws.Cells.Item(19, 14).Value = 1.366417278529475
$ws.Cells.Item(19, 15).Value = 3.507873628741436
$ws.Cells.Item(20, 2).Value = 1.358279025233287
$ws.Cells.Item(20, 3).Value = 0.2476567286588534
$ws.Cells.Item(20, 4).Value = 0.1043555544411419
$ws.Cells.Item(20, 5).Value = 0.08762460424203322
$ws.Cells.Item(20, 7).Value = 0.00243921819108626
$ws.Cells.Item(20, 9).Value = 0.8011763536673797
$ws.Cells.Item(20, 12).Value = 0.2183879643218489
$ws.Cells.Item(20, 13).Value = 0.2874441547633992
$ws.Cells.Item(20, 14).Value = 1.361054940466403
$ws.Cells.Item(20, 15).Value = 3.516647164373921
$ws.Cells.Item(21, 2).Value = 1.481996760586355
$ws.Cells.Item(21, 3).Value = 0.2715017316089927
$ws.Cells.Item(21, 4).Value = 0.1149018535734854
$ws.Cells.Item(21, 5).Value = 0.08698493244836869
$ws.Cells.Item(21, 7).Value = 0.002435570390548425
$ws.Cells.Item(21, 9).Value = 0.7980862946605924
$ws.Cells.Item(21, 12).Value = 0.2231200599594132
$ws.Cells.Item(21, 13).Value = 0.3074609012610665
$ws.Cells.Item(21, 14).Value = 1.34367268248781
$ws.Cells.Item(21, 15).Value = 3.548789854428605
$ws.Cells.Item(22, 2).Value = 1.563052139216666
$ws.Cells.Item(22, 3).Value = 0.2870153786429057
$ws.Cells.Item(22, 4).Value = 0.1218283987212203
$ws.Cells.Item(22, 5).Value = 0.08658453946455791
$ws.Cells.Item(22, 7).Value = 0.002433277533983916
$ws.Cells.Item(22, 9).Value = 0.7964887062997832
$ws.Cells.Item(22, 12).Value = 0.2263035330294088
$ws.Cells.Item(22, 13).Value = 0.320617556170177
$ws.Cells.Item(22, 14).Value = 1.332785572852199
$ws.Cells.Item(22, 15).Value = 3.571787161351921
$ws.Cells.Item(23, 2).Value = 1.519773558379995
$ws.Cells.Item(23, 3).Value = 0.2787419350758
$ws.Cells.Item(23, 4).Value = 0.118128457850105
$ws.Cells.Item(23, 5).Value = 0.08679662253944764
$ws.Cells.Item(23, 7).Value = 0.002434492955247087
$ws.Cells.Item(23, 9).Value = 0.7973027126632886
$ws.Cells.Item(23, 12).Value = 0.2245961614798233
$ws.Cells.Item(23, 13).Value = 0.3135888394036712
$ws.Cells.Item(23, 14).Value = 1.338552908570762
$ws.Cells.Item(23, 15).Value = 3.559330700151691
$ws.Cells.Item(24, 2).Value = 1.356260362059516
$ws.Cells.Item(24, 3).Value = 0.2472658133921186
$ws.Cells.Item(24, 4).Value = 0.1041837514435855
$ws.Cells.Item(24, 5).Value = 0.08763535821906221
$ws.Cells.Item(24, 7).Value = 0.002439279357472234
$ws.Cells.Item(24, 9).Value = 0.8012339715060079
$ws.Cells.Item(24, 12).Value = 0.2183121562369763
$ws.Cells.Item(24, 13).Value = 0.2871182632700808
$ws.Cells.Item(24, 14).Value = 1.361347028375889
$ws.Cells.Item(24, 15).Value = 3.516155348969932
$ws.Cells.Item(25, 2).Value = 1.180775734696056
$ws.Cells.Item(25, 3).Value = 0.2130148592328567
$ws.Cells.Item(25, 4).Value = 0.08928683131432535
$ws.Cells.Item(25, 5).Value = 0.08861619500182483
$ws.Cells.Item(25, 7).Value = 0.002444836311374912
$ws.Cells.Item(25, 9).Value = 0.8072811119820358
$ws.Cells.Item(25, 12).Value = 0.211924940968423
$ws.Cells.Item(25, 13).Value = 0.2588914336270989
$ws.Cells.Item(25, 14).Value = 1.387963803326919
$ws.Cells.Item(25, 15).Value = 3.47810166147579
